$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Counterspell', ['{U}{U}', 'Instant', 'Counter target spell.'])"
$ws.Range("A3").Value = "('Vampiric Tutor', ['{B}', 'Instant', 'Search your library for a card, then shuffle your library and put that card on top of it. You lose 2 life.'])"
$ws.Range("A4:A9").ClearContents()
